# budget_target_2425.xlsx update
# - Row 3 (Arsenal Women, UEFA Women's Champions League budget line):
#   fixture renamed from "Arsenal Women UWCL Quarter-Final (Date TBC)"
#   to "Arsenal Women v Real Madrid Women" (the fixture became confirmed).
# - Row 6 (Arsenal, UEFA Champions League budget line):
#   the "West Ham United" fixture's budget target is revised
#   from 490113 to 712500.
# - Scroll the sheet view so row 4 is at the top (cosmetic, matches the
#   author's saved window position) without disturbing the existing
#   C7 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Arsenal Women v Real Madrid Women"
$ws.Range("C6").Value = 712500

$aw = $excel.ActiveWindow()
$aw.ScrollRow = 4
